$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.7888409
$ws.Range("K2").Value = 12
$ws.Range("M2").Value = 0.9178586
$ws.Range("O2").Value = 12
$ws.Range("Q2").Value = 0.9430005
$ws.Range("U2").Value = 0.9374105
$ws.Range("Y2").Value = 0.9341104
$ws.Range("AC2").Value = 9.8
$ws.Range("I3").Value = 0.7936414
$ws.Range("M3").Value = 0.9209819
$ws.Range("Q3").Value = 0.9427314
$ws.Range("U3").Value = 0.939803
$ws.Range("W3").Value = 3
$ws.Range("Y3").Value = 0.9194454
$ws.Range("AA3").Value = 18
$ws.Range("AC3").Value = 11.6
$ws.Range("I4").Value = 0.7961511
$ws.Range("M4").Value = 0.9225368
$ws.Range("O4").Value = 6
$ws.Range("Q4").Value = 0.9444176
$ws.Range("U4").Value = 0.9401447
$ws.Range("W4").Value = 2
$ws.Range("Y4").Value = 0.9203255
$ws.Range("AA4").Value = 16
$ws.Range("AC4").Value = 7.8
$ws.Range("I5").Value = 0.7885364
$ws.Range("K5").Value = 15
$ws.Range("M5").Value = 0.9159912
$ws.Range("O5").Value = 15
$ws.Range("Q5").Value = 0.9428095
$ws.Range("S5").Value = 14
$ws.Range("U5").Value = 0.9320369
$ws.Range("Y5").Value = 0.9293631
$ws.Range("AA5").Value = 11
$ws.Range("AC5").Value = 14.2
$ws.Range("I6").Value = 0.7885344
$ws.Range("K6").Value = 16
$ws.Range("M6").Value = 0.915987
$ws.Range("O6").Value = 16
$ws.Range("Q6").Value = 0.9428095
$ws.Range("U6").Value = 0.9320369
$ws.Range("Y6").Value = 0.9293631
$ws.Range("AA6").Value = 12
$ws.Range("AC6").Value = 15.2
$ws.Range("I7").Value = 0.7885344
$ws.Range("K7").Value = 17
$ws.Range("M7").Value = 0.915987
$ws.Range("O7").Value = 17
$ws.Range("Q7").Value = 0.9428095
$ws.Range("U7").Value = 0.9320369
$ws.Range("Y7").Value = 0.9293631
$ws.Range("AA7").Value = 13
$ws.Range("AC7").Value = 16.2
$ws.Range("I8").Value = 0.7888388
$ws.Range("K8").Value = 14
$ws.Range("M8").Value = 0.9178557
$ws.Range("O8").Value = 14
$ws.Range("Q8").Value = 0.9430005
$ws.Range("S8").Value = 10
$ws.Range("U8").Value = 0.9374105
$ws.Range("Y8").Value = 0.9341104
$ws.Range("AC8").Value = 11.2
$ws.Range("I9").Value = 0.7936414
$ws.Range("M9").Value = 0.9209819
$ws.Range("Q9").Value = 0.9427301
$ws.Range("S9").Value = 20
$ws.Range("U9").Value = 0.9396098
$ws.Range("W9").Value = 5
$ws.Range("Y9").Value = 0.9194454
$ws.Range("AA9").Value = 19
$ws.Range("AC9").Value = 12.8
$ws.Range("I10").Value = 0.7959714
$ws.Range("M10").Value = 0.9225368
$ws.Range("Q10").Value = 0.9444582
$ws.Range("U10").Value = 0.9397074
$ws.Range("W10").Value = 4
$ws.Range("Y10").Value = 0.9203255
$ws.Range("AA10").Value = 17
$ws.Range("AC10").Value = 8.6
$ws.Range("I11").Value = 0.7888409
$ws.Range("K11").Value = 13
$ws.Range("M11").Value = 0.9178574
$ws.Range("O11").Value = 13
$ws.Range("Q11").Value = 0.9430005
$ws.Range("S11").Value = 11
$ws.Range("U11").Value = 0.9374105
$ws.Range("Y11").Value = 0.9341104
$ws.Range("AC11").Value = 11.4
$ws.Range("I12").Value = 0.7936414
$ws.Range("M12").Value = 0.9209819
$ws.Range("Q12").Value = 0.9427623
$ws.Range("S12").Value = 17
$ws.Range("U12").Value = 0.9396098
$ws.Range("W12").Value = 6
$ws.Range("Y12").Value = 0.9194454
$ws.Range("AA12").Value = 20
$ws.Range("AC12").Value = 13
$ws.Range("I13").Value = 0.7885344
$ws.Range("K13").Value = 18
$ws.Range("M13").Value = 0.9159293
$ws.Range("O13").Value = 20
$ws.Range("Q13").Value = 0.942742
$ws.Range("S13").Value = 18
$ws.Range("U13").Value = 0.9414732
$ws.Range("Y13").Value = 0.931004
$ws.Range("AC13").Value = 13.4
$ws.Range("I14").Value = 0.8015114
$ws.Range("M14").Value = 0.9246578
$ws.Range("O14").Value = 2
$ws.Range("Q14").Value = 0.9470992
$ws.Range("U14").Value = 0.937812
$ws.Range("Y14").Value = 0.9341104
$ws.Range("AC14").Value = 4.2
$ws.Range("I15").Value = 0.7884538
$ws.Range("K15").Value = 20
$ws.Range("M15").Value = 0.9159723
$ws.Range("O15").Value = 19
$ws.Range("Q15").Value = 0.9428186
$ws.Range("S15").Value = 13
$ws.Range("U15").Value = 0.9320369
$ws.Range("Y15").Value = 0.9293527
$ws.Range("AA15").Value = 14
$ws.Range("AC15").Value = 17
$ws.Range("I16").Value = 0.8005746
$ws.Range("M16").Value = 0.9246582
$ws.Range("Q16").Value = 0.9471593
$ws.Range("U16").Value = 0.937812
$ws.Range("Y16").Value = 0.9341104
$ws.Range("I17").Value = 0.7964122
$ws.Range("M17").Value = 0.923434
$ws.Range("Q17").Value = 0.9465523
$ws.Range("U17").Value = 0.9380812
$ws.Range("Y17").Value = 0.9343387
$ws.Range("I18").Value = 0.8016446
$ws.Range("M18").Value = 0.9246262
$ws.Range("Q18").Value = 0.9468049
$ws.Range("U18").Value = 0.9376773
$ws.Range("Y18").Value = 0.9341104
$ws.Range("I19").Value = 0.7884559
$ws.Range("K19").Value = 19
$ws.Range("M19").Value = 0.9159727
$ws.Range("O19").Value = 18
$ws.Range("Q19").Value = 0.9428359
$ws.Range("U19").Value = 0.9320369
$ws.Range("Y19").Value = 0.9293527
$ws.Range("AA19").Value = 15
$ws.Range("AC19").Value = 16.8
$ws.Range("I20").Value = 0.800724
$ws.Range("M20").Value = 0.9246567
$ws.Range("O20").Value = 3
$ws.Range("Q20").Value = 0.9468523
$ws.Range("U20").Value = 0.9376773
$ws.Range("Y20").Value = 0.9341104
$ws.Range("AC20").Value = 6
$ws.Range("I21").Value = 0.7964637
$ws.Range("M21").Value = 0.922374
$ws.Range("O21").Value = 8
$ws.Range("Q21").Value = 0.945688
$ws.Range("U21").Value = 0.9379466
$ws.Range("Y21").Value = 0.9344553
$ws.Range("AC21").Value = 5.6
